# Denormalize db schema: remove many-to-many relation table by adding
# "tag_ids" and "doc_ids" columns directly onto the folder table.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- 1. Grow the Excel Table (ListObject) from 16 to 18 columns -----------
$lo = $ws.ListObjects.Item(1)
$lo.Resize($ws.Range("A1:R26"))

# --- 2. Header row for the two new columns (this also names the table
#        columns, since the table reads the header cells on save) ----------
$ws.Range("Q1").Value = "tag_ids"
$ws.Range("R1").Value = "doc_ids"

# --- 3. New column data -----------------------------------------------------
# row 13 = "10-tourisme" folder gains a doc_ids example
$ws.Range("R13").Value = "tourisme_exemple"

# row 26 = "bevnat" folder gains doc_ids / tag_ids
$ws.Range("R26").Value = "bevnat_info, bevnat_variable"

# row 25 = "statpop" folder gains doc_ids / tag_ids
$ws.Range("R25").Value = "statpop_info"

# row 26 = "bevnat" folder tag_ids
$ws.Range("Q26").Value = "personal_data, sensible_data, population"

# row 25 = "statpop" folder tag_ids
$ws.Range("Q25").Value = "personal_data, population, societe"

# --- 4. Column widths for the two new columns ------------------------------
$ws.Columns.Item(17).ColumnWidth = 32
$ws.Columns.Item(18).ColumnWidth = 22.5

# --- 5. Freeze both header row and first column, restore selection --------
$win = $excel.ActiveWindow
$win.FreezePanes = $false
[void]$ws.Range("B2").Select()
$win.FreezePanes = $true
[void]$ws.Range("U25").Select()
